# Update "想去人数" (F column) counts, and for one row the "最低票价" (G
# column) status, across the "展览" (sheet 1), "演出" (sheet 2) and
# "全部类型" (sheet 4, the combined roll-up) worksheets — matching the
# refreshed scrape output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) -------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = 2517
$ws1.Range("F3").Value = 339
$ws1.Range("G3").Value = 218      # was inline text "已售罄" (sold out); now has a real min price
$ws1.Range("F4").Value = 346
$ws1.Range("F5").Value = 1426
$ws1.Range("F6").Value = 1112
$ws1.Range("F7").Value = 317
$ws1.Range("F8").Value = 525
$ws1.Range("F11").Value = 105
$ws1.Range("F12").Value = 549
$ws1.Range("F13").Value = 8772
$ws1.Range("F14").Value = 378
$ws1.Range("F16").Value = 266
$ws1.Range("F18").Value = 177
$ws1.Range("F20").Value = 601
$ws1.Range("F22").Value = 1161
$ws1.Range("F23").Value = 1001
$ws1.Range("F24").Value = 2044
$ws1.Range("F25").Value = 2109
$ws1.Range("F27").Value = 1800
$ws1.Range("F29").Value = 1918
$ws1.Range("F31").Value = 222
$ws1.Range("F33").Value = 111
$ws1.Range("F34").Value = 193
$ws1.Range("F35").Value = 12
$ws1.Range("F36").Value = 312
$ws1.Range("F38").Value = 260
$ws1.Range("F39").Value = 444
$ws1.Range("F40").Value = 754
$ws1.Range("F42").Value = 270

# --- Sheet 2: 演出 (Performances) ------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("F4").Value = 9

# --- Sheet 4: 全部类型 (All types, combined roll-up) ------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value = 2517
$ws4.Range("F3").Value = 339
$ws4.Range("G3").Value = 218      # was inline text "已售罄" (sold out); now has a real min price
$ws4.Range("F4").Value = 346
$ws4.Range("F5").Value = 1426
$ws4.Range("F7").Value = 1112
$ws4.Range("F8").Value = 317
$ws4.Range("F9").Value = 525
$ws4.Range("F12").Value = 105
$ws4.Range("F13").Value = 549
$ws4.Range("F14").Value = 8772
$ws4.Range("F15").Value = 378
$ws4.Range("F18").Value = 266
$ws4.Range("F20").Value = 177
$ws4.Range("F22").Value = 601
$ws4.Range("F24").Value = 1161
$ws4.Range("F25").Value = 1001
$ws4.Range("F26").Value = 2044
$ws4.Range("F27").Value = 2109
$ws4.Range("F29").Value = 1800
$ws4.Range("F31").Value = 1918
$ws4.Range("F33").Value = 222
$ws4.Range("F35").Value = 111
$ws4.Range("F36").Value = 193
$ws4.Range("F37").Value = 12
$ws4.Range("F38").Value = 312
$ws4.Range("F40").Value = 260
$ws4.Range("F41").Value = 444
$ws4.Range("F42").Value = 9
$ws4.Range("F46").Value = 754
$ws4.Range("F49").Value = 270
